$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain text (coin names, links, percentages,
# or numbers that already contain thousands separators so Excel will not
# reinterpret them as numeric values).
$textUpdates = @{
    "D2" = "70.033.84"
    "E2" = "  -0.64%  "
    "D3" = "3.476.09"
    "E3" = "  -1.92%  "
    "E4" = "  -0.32%  "
    "E5" = "  +1.77%  "
    "E6" = "  -3.19%  "
    "D7" = "3.475.07"
    "E7" = "  -1.77%  "
    "E8" = "  -2.40%  "
    "E9" = "  -0.11%  "
    "E10" = "  -0.47%  "
    "E11" = "  -1.79%  "
    "E12" = "  -3.35%  "
    "E13" = "  -3.88%  "
    "E14" = "  -2.57%  "
    "D15" = "4.038.07"
    "E15" = "  -2.28%  "
    "E16" = "  -1.20%  "
    "E17" = "  -3.01%  "
    "B18" = "WrappedBTC"
    "C18" = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
    "D18" = "70.137.34"
    "E18" = "  -0.80%  "
    "B19" = "WrappedEther"
    "C19" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "D19" = "3.471.03"
    "E19" = "  -2.47%  "
    "E20" = "  +0.77%  "
    "E21" = "  -0.28%  "
    "E22" = "  -2.10%  "
    "E23" = "  -4.38%  "
    "E24" = "  -0.85%  "
    "E25" = "  -2.88%  "
    "E26" = "  -2.82%  "
    "E27" = "  -0.05%  "
    "E28" = "  -4.62%  "
    "E29" = "  -2.61%  "
    "E30" = "  -4.08%  "
    "E31" = "  -4.01%  "
    "E32" = "  -6.72%  "
    "E33" = "  -3.24%  "
    "E34" = "  -6.10%  "
    "E35" = "  -21.02%  "
    "E36" = "  -0.80%  "
    "E37" = "  +0.68%  "
    "E38" = "  -4.08%  "
    "E39" = "  +0.01%  "
    "E40" = "  -0.78%  "
    "E41" = "  -0.89%  "
    "E42" = "  -10.95%  "
    "D43" = "3.262.53"
    "E43" = "  -2.92%  "
    "D44" = "0.0₃0699"
    "E44" = "  +0.37%  "
    "E45" = "  -6.31%  "
    "E46" = "  -5.08%  "
    "E47" = "  -5.79%  "
    "E48" = "  -6.64%  "
    "E49" = "  -2.62%  "
    "E50" = "  -0.09%  "
}

# Cells whose new values look like plain numbers (e.g. "613.84"). The source
# workbook stores these as text (inline strings), so we force the cell to
# text format before assigning the value, then drop the now-unneeded
# number format so the cell keeps its original (default) styling.
$numericLookingUpdates = @{
    "D4" = "0.999"
    "D5" = "613.84"
    "D6" = "167.93"
    "D10" = "0.193"
    "D11" = "7.13"
    "D12" = "0.566"
    "D13" = "44.56"
    "D16" = "8.21"
    "D17" = "589.32"
    "D21" = "17.24"
    "D22" = "0.856"
    "D23" = "8.87"
    "D24" = "95.65"
    "D25" = "15.22"
    "D28" = "2.48"
    "D29" = "33.10"
    "D30" = "8.67"
    "D31" = "7.86"
    "D34" = "6.59"
    "D35" = "580.91"
    "D36" = "10.66"
    "D37" = "0.0482"
    "D38" = "0.0963"
    "D40" = "56.32"
    "D41" = "0.141"
    "D42" = "3.18"
    "D45" = "0.296"
    "D46" = "30.85"
    "D50" = "133.81"
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

foreach ($ref in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$ref]
    $cell.ClearFormats()
}
